$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) contain duplicate rows.
# Row 7 (event "合肥·第十三届次元之门动漫游戏博览会"): F7 6974 -> 6975
# Row 16 (event "合肥·首届佀活企划——佀像计划-阎罗舞台"): F16 22 -> 23
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F7").Value = 6975
    $ws.Range("F16").Value = 23
}
